# Weekly refresh of the "Cilantro" price series: a new weekly record is
# inserted at row 117 (pushing the existing history for that market down by
# one row), and the worksheet dimension grows from A1:R173 to A1:R174.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row above the current row 117, shifting rows 117:173 down to
# 118:174 (this also carries the date-format style of row 117 onto the new
# row, matching column D's existing formatting).
$ws.Rows("117:117").Insert()

# Populate the newly inserted row with this week's record.
$ws.Range("A117").Value = 8
$ws.Range("B117").Value = "Terminal La Palmera de La Serena"
$ws.Range("C117").Value = "Coquimbo"
$ws.Range("D117").Value = 44845
$ws.Range("E117").Value = 4
$ws.Range("F117").Value = 100112040
$ws.Range("G117").Value = "Cilantro"
$ws.Range("H117").Value = "Sin especificar"
$ws.Range("I117").Value = "Primera"
$ws.Range("J117").Value = 3200
$ws.Range("K117").Value = 2000
$ws.Range("L117").Value = 2500
$ws.Range("M117").Value = 2250
$ws.Range("N117").Value = "`$/atado 1 a 1,5 kilos"
$ws.Range("O117").Value = "Provincia del Elquí"
$ws.Range("P117").Value = 1500
$ws.Range("Q117").Value = 1.5
$ws.Range("R117").Value = "Hortaliza"
